# Apply updated PvsI model-fitting results to final rates sheet (500_rates.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1405756097560976
$ws.Range("Z2").Value = -0.117334580617608
$ws.Range("AB2").Value = -473.4144617474702
$ws.Range("AD2").Value = -473.4144617474702

# Row 3
$ws.Range("T3").Value = 0.144009756097561
$ws.Range("Z3").Value = -0.1731737476119215
$ws.Range("AB3").Value = -1140.826853027863
$ws.Range("AD3").Value = -1140.826853027863

# Row 4
$ws.Range("T4").Value = 0.1462634146341464
$ws.Range("Z4").Value = -0.1315190703102171
$ws.Range("AB4").Value = -842.5214905140799
$ws.Range("AD4").Value = -842.5214905140799

# Row 5
$ws.Range("T5").Value = 0.1500390243902439
$ws.Range("Z5").Value = -0.1270169707930712
$ws.Range("AB5").Value = -793.2549583356575
$ws.Range("AD5").Value = -793.2549583356575

# Row 6
$ws.Range("T6").Value = 0.1463707317073171
$ws.Range("Z6").Value = -0.1000315627773721
$ws.Range("AB6").Value = -624.8953060000269
$ws.Range("AD6").Value = -624.8953060000269

# Row 7
$ws.Range("T7").Value = 0.1465658536585366
$ws.Range("Z7").Value = -0.1653736955078527
$ws.Range("AB7").Value = -825.200965484603
$ws.Range("AD7").Value = -825.200965484603

# Row 8
$ws.Range("T8").Value = 0.1449658536585366
$ws.Range("Z8").Value = -0.1666309957274192
$ws.Range("AB8").Value = -672.133957072556
$ws.Range("AD8").Value = -672.133957072556

# Row 9
$ws.Range("T9").Value = 0.1544
$ws.Range("Z9").Value = -0.001356116654488074

# Row 10
$ws.Range("T10").Value = 0.1405756097560976
$ws.Range("Z10").Value = 0.1400890362927661
$ws.Range("AB10").Value = 565.2227618164701
$ws.Range("AD10").Value = 565.2227618164701

# Row 11
$ws.Range("T11").Value = 0.144009756097561
$ws.Range("Z11").Value = 0.1163457900251444
$ws.Range("AB11").Value = 766.4579841216574
$ws.Range("AD11").Value = 766.4579841216574

# Row 12
$ws.Range("T12").Value = 0.1462634146341464
$ws.Range("Z12").Value = 0.05588008741716961
$ws.Range("AB12").Value = 357.9722273714504
$ws.Range("AD12").Value = 357.9722273714504

# Row 13
$ws.Range("T13").Value = 0.1500390243902439
$ws.Range("Z13").Value = 0.119309848809772
$ws.Range("AB13").Value = 745.1219199741149
$ws.Range("AD13").Value = 745.1219199741149

# Row 14
$ws.Range("T14").Value = 0.1463707317073171
$ws.Range("Z14").Value = 0.04624812368653618
$ws.Range("AB14").Value = 288.9111656422364
$ws.Range("AD14").Value = 288.9111656422364

# Row 15
$ws.Range("T15").Value = 0.1465658536585366
$ws.Range("Z15").Value = 0.161351139913875
$ws.Range("AB15").Value = 805.1287481366617
$ws.Range("AD15").Value = 805.1287481366617

# Row 16
$ws.Range("T16").Value = 0.1449658536585366
$ws.Range("Z16").Value = 0.2626824278868477
$ws.Range("AB16").Value = 1059.573454135944
$ws.Range("AD16").Value = 1059.573454135944

# Row 17
$ws.Range("T17").Value = 0.1544
$ws.Range("Z17").Value = 0.005649580567729108

